$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Split the final paragraph ("...In sass, & take the place of the parent
#    's name." + _GoBack bookmark) into three paragraphs:
#      a) the original sentence (unchanged, bookmark removed from here)
#      b) "Update set-value and meme in package-lock.json: use "
#      c) "npm install set-value and npm install mem" + _GoBack bookmark +
#         a trailing space run
# ---------------------------------------------------------------------------

$lastPara = $d.Paragraphs.Last
$splitPos = $lastPara.Range.End - 1   # right after "name." / right before the bookmark

# The _GoBack bookmark currently sits at $splitPos; remove it here, it gets
# re-created (in its new home) by the XML inserted below.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$insertRange = $d.Range($splitPos, $splitPos)

$newParasXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="6"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
              <w:ind w:firstLineChars="0"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="default"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve">Update set-value and meme in package-lock.json: use </w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="6"/>
              <w:numPr>
                <w:numId w:val="0"/>
              </w:numPr>
              <w:ind w:leftChars="0"/>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="default"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>npm install set-value and npm install mem</w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="default"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertRange.InsertXML($newParasXml)

# ---------------------------------------------------------------------------
# 2) Mark the "FollowedHyperlink" style as a Quick Style (w:qFormat) - this
#    is what Word does when a style is promoted to the Quick Style gallery.
# ---------------------------------------------------------------------------

$followedHyperlink = $d.Styles("FollowedHyperlink")
$followedHyperlink.QuickStyle = $true

Write-Output "done"
